$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in computed metrics for rows 2-7 (Clinical only / Clinical + wb-FA)
$ws.Range("C2").Value = 0.639
$ws.Range("D2").Value = 14.741
$ws.Range("E2").Value = 11.547
$ws.Range("F2").Value = 235.454
$ws.Range("G2").Value = 6

$ws.Range("C3").Value = -1.593
$ws.Range("D3").Value = 34.188
$ws.Range("E3").Value = 25.427
$ws.Range("F3").Value = 235.418
$ws.Range("G3").Value = 6

$ws.Range("C4").Value = -3.08
$ws.Range("D4").Value = 50.779
$ws.Range("E4").Value = 36.469
$ws.Range("F4").Value = 236.177
$ws.Range("G4").Value = 6

$ws.Range("C5").Value = 0.651
$ws.Range("D5").Value = 14.497
$ws.Range("E5").Value = 11.273
$ws.Range("F5").Value = 237.452
$ws.Range("G5").Value = 6

$ws.Range("C6").Value = -2.14
$ws.Range("D6").Value = 37.621
$ws.Range("E6").Value = 29.322
$ws.Range("F6").Value = 234.607
$ws.Range("G6").Value = 6

$ws.Range("C7").Value = -3.432
$ws.Range("D7").Value = 52.92
$ws.Range("E7").Value = 37.818
$ws.Range("F7").Value = 237.84
$ws.Range("G7").Value = 6

# Clear stray zero values left over in G8 and G10 (Clinical + tractFA rows)
$ws.Range("G8").Value = ""
$ws.Range("G10").Value = ""
